$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "FECHA REAL DE CIERRE" (actual closing date) values that were
# missing for rows 4-8 (the dates mirror the "planned" dates in column D).
$ws.Range("E4").Value = 42361
$ws.Range("E5").Value = 42361
$ws.Range("E6").Value = 42359
$ws.Range("E7").Value = 42359
$ws.Range("E8").Value = 42359

# Row 8 shrinks slightly once the text re-flows around the new date.
$ws.Rows.Item(8).RowHeight = 28.35

# The closing-date cell for row 8 picks up a distinct (but visually
# identical) font as part of the correction.
$ws.Range("E8").Font.ThemeFont = 1

# Leave the selection where the editor left it.
[void]$ws.Range("E7").Select()
